$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 463, pushing the existing rows 463:506
# down to 465:508 (new weekly price records added at the top of this
# date-sorted block).
$ws.Range("463:464").Insert()

# New row 463: Camote, "1a (guarda)"
$ws.Range("A463").Value = 5
$ws.Range("B463").Value = "Macroferia Regional de Talca"
$ws.Range("C463").Value = "Maule"
$ws.Range("D463").Value = 45194
$ws.Range("E463").Value = 7
$ws.Range("F463").Value = 100112045
$ws.Range("G463").Value = "Zapallo"
$ws.Range("H463").Value = "Camote"
$ws.Range("I463").Value = "1a (guarda)"
$ws.Range("J463").Value = 400
$ws.Range("K463").Value = 650
$ws.Range("L463").Value = 650
$ws.Range("M463").Value = 650
$ws.Range("N463").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O463").Value = "Región del Maule"
$ws.Range("P463").Value = 650
$ws.Range("Q463").Value = 1
$ws.Range("R463").Value = "Hortaliza"

# New row 464: Paine, "1a (guarda)"
$ws.Range("A464").Value = 5
$ws.Range("B464").Value = "Macroferia Regional de Talca"
$ws.Range("C464").Value = "Maule"
$ws.Range("D464").Value = 45194
$ws.Range("E464").Value = 7
$ws.Range("F464").Value = 100112045
$ws.Range("G464").Value = "Zapallo"
$ws.Range("H464").Value = "Paine"
$ws.Range("I464").Value = "1a (guarda)"
$ws.Range("J464").Value = 1500
$ws.Range("K464").Value = 400
$ws.Range("L464").Value = 400
$ws.Range("M464").Value = 400
$ws.Range("N464").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O464").Value = "Región del Maule"
$ws.Range("P464").Value = 400
$ws.Range("Q464").Value = 1
$ws.Range("R464").Value = "Hortaliza"
